# Rename the "study" worksheet to "defaultValues" and make it the active/
# selected sheet (tab) instead of "OTU_table".
#
# Context (from commit message "allow user to enter dataset title on
# upload, changes to users dataset list"): the former "study" template
# sheet is repurposed as a "defaultValues" sheet, and the workbook is saved
# with that sheet active/visible (its tab selected) in place of the first
# sheet.

$wb = $excel.ActiveWorkbook

# Sheet 4 ("study") -> rename to "defaultValues"
$wsDefaults = $wb.Worksheets.Item(4)
$wsDefaults.Name = "defaultValues"

# Make the renamed sheet the active tab (this both updates the workbook's
# activeTab and moves tabSelected from the previously active sheet to this
# one).
$wsDefaults.Activate()
